$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "TC22018"

# Populate the new "Mounted?"-style Y/Y columns (I & J) for rows 287-306
for ($r = 287; $r -le 306; $r++) {
    $ws.Cells.Item($r, 9).Value = "Y"
    $ws.Cells.Item($r, 10).Value = "Y"
}

# Update the view state: zoom + scroll position + selection to match the saved view
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollRow = 279
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I287:J306").Select()
